# Applies the cell updates described by the commit diff to the active workbook.
# Columns B/C (coin name / link) are plain text; columns D/E (price / volume%)
# look numeric, so they are written with a leading apostrophe to force Excel to
# store them as literal text (preserving exact formatting such as trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'314.54"
$ws.Range("E2").Value = "'2.31%"
$ws.Range("D3").Value = "'40.90"
$ws.Range("E3").Value = "'-0.24%"
$ws.Range("D4").Value = "'5.167"
$ws.Range("E4").Value = "'-1.24%"
$ws.Range("D5").Value = "'0.07603"
$ws.Range("E5").Value = "'-0.74%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.682"
$ws.Range("E6").Value = "'2.90%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9311"
$ws.Range("E7").Value = "'1.77%"
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D8").Value = "'0.1202"
$ws.Range("E8").Value = "'-3.34%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1824"
$ws.Range("E9").Value = "'0.07%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.09058"
$ws.Range("E10").Value = "'-0.88%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.04136"
$ws.Range("E11").Value = "'-0.56%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.1055"
$ws.Range("E12").Value = "'0.40%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001281"
$ws.Range("E13").Value = "'1.44%"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").Value = "'0.005902"
$ws.Range("E14").Value = "'1.21%"
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").Value = "'3.336"
$ws.Range("E15").Value = "'-0.24%"
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").Value = "'4.320"
$ws.Range("E16").Value = "'0.21%"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").Value = "'2.424"
$ws.Range("E17").Value = "'-0.53%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "'0.3348"
$ws.Range("E18").Value = "'0.39%"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "'7.597"
$ws.Range("E19").Value = "'1.59%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1351"
$ws.Range("E20").Value = "'-3.01%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "'0.2832"
$ws.Range("E21").Value = "'-1.79%"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "'0.03977"
$ws.Range("E22").Value = "'-2.08%"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "'0.001279"
$ws.Range("E23").Value = "'1.30%"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "'0.004079"
$ws.Range("E24").Value = "'-4.67%"
$ws.Range("B25").Value = "NitroEx"
$ws.Range("C25").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D25").Value = "'0.0001349"
$ws.Range("E25").Value = "'5.98%"
$ws.Range("B26").Value = "UpBots"
$ws.Range("C26").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D26").Value = "'0.0003039"
$ws.Range("E26").Value = "'-95.95%"
$ws.Range("E38").Value = "'-2.31%"
$ws.Range("D39").Value = "'0.05183"
$ws.Range("E39").Value = "'-2.89%"
$ws.Range("D40").Value = "'0.007699"
$ws.Range("E40").Value = "'-2.03%"
$ws.Range("D41").Value = "'0.1303"
$ws.Range("E41").Value = "'-0.75%"
$ws.Range("D42").Value = "'0.007588"
$ws.Range("E42").Value = "'13.85%"
$ws.Range("D43").Value = "'0.003300"
$ws.Range("E43").Value = "'72.35%"
$ws.Range("D44").Value = "'0.008481"
$ws.Range("E44").Value = "'10.44%"
$ws.Range("D45").Value = "'0.3394"
$ws.Range("E45").Value = "'10.99%"
$ws.Range("D46").Value = "'0.00006589"
$ws.Range("E46").Value = "'-1.82%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.30%"
$ws.Range("D48").Value = "'0.2737"
$ws.Range("E48").Value = "'-36.69%"
$ws.Range("D49").Value = "'0.004200"
$ws.Range("E49").Value = "'35.15%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.30%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.30%"
